$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "Une page produit affichant de manière dynamique les détails du produit sur lequel l'utilisateur a cliqué"
$ws.Range("B4").Value = "Possibilité sur les pages produits de sélectionner ou entrer sa quantité, changer la couleur du canapé et l'ajouter au panier"

# Row heights reflow after the wrapped text grows rows 3 & 4 (and row 2
# picks up the same auto-fit height as its sibling header-data rows).
$ws.Rows.Item(2).RowHeight = 54
$ws.Rows.Item(3).RowHeight = 54
$ws.Rows.Item(4).RowHeight = 54
$ws.Rows.Item(7).RowHeight = 18
for ($r = 5; $r -le 22; $r++) {
    if ($r -ne 7) {
        $ws.Rows.Item($r).RowHeight = 17
    }
}

$null = $ws.Range("B5").Select()
